# Fruta / hortaliza, semanal
# Insert two new price-record rows (Primera / Segunda) for a new survey
# date (2022-12-23, serial 44918) right before the existing row 812 block,
# pushing the rest of the "Brócoli" records down by two rows
# (old 812..909 -> new 814..911).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 812 downward (two rows) to make room for the new records.
$ins = $ws.Range("A812:R813")
$ins.Insert()

# New row 812 - "Primera" quality record for 2022-12-23
$ws.Cells.Item(812,1).Value  = 3
$ws.Cells.Item(812,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(812,3).Value  = "Coquimbo"
$ws.Cells.Item(812,4).Value  = 44918
$ws.Cells.Item(812,5).Value  = 5
$ws.Cells.Item(812,6).Value  = 100112023
$ws.Cells.Item(812,7).Value  = "Brócoli"
$ws.Cells.Item(812,8).Value  = "Sin especificar"
$ws.Cells.Item(812,9).Value  = "Primera"
$ws.Cells.Item(812,10).Value = 2300
$ws.Cells.Item(812,11).Value = 800
$ws.Cells.Item(812,12).Value = 850
$ws.Cells.Item(812,13).Value = 824
$ws.Cells.Item(812,14).Value = "$/unidad"
$ws.Cells.Item(812,15).Value = "Provincia de Quillota"
$ws.Cells.Item(812,16).Value = 824
$ws.Cells.Item(812,17).Value = 1
$ws.Cells.Item(812,18).Value = "Hortaliza"

# New row 813 - "Segunda" quality record for 2022-12-23
$ws.Cells.Item(813,1).Value  = 3
$ws.Cells.Item(813,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(813,3).Value  = "Coquimbo"
$ws.Cells.Item(813,4).Value  = 44918
$ws.Cells.Item(813,5).Value  = 5
$ws.Cells.Item(813,6).Value  = 100112023
$ws.Cells.Item(813,7).Value  = "Brócoli"
$ws.Cells.Item(813,8).Value  = "Sin especificar"
$ws.Cells.Item(813,9).Value  = "Segunda"
$ws.Cells.Item(813,10).Value = 1200
$ws.Cells.Item(813,11).Value = 700
$ws.Cells.Item(813,12).Value = 700
$ws.Cells.Item(813,13).Value = 700
$ws.Cells.Item(813,14).Value = "$/unidad"
$ws.Cells.Item(813,15).Value = "Provincia de Quillota"
$ws.Cells.Item(813,16).Value = 700
$ws.Cells.Item(813,17).Value = 1
$ws.Cells.Item(813,18).Value = "Hortaliza"
